$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F holds "dSF". Update specific rows per the repull/recalculation.
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -2
$ws.Range("F12").Value = -6
